$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "estado" -> "Estado"
$ws.Range("I1").Value = "Estado"

# Update existing row 2 values
$ws.Range("D2").Value = 1200.5446543232
$ws.Range("F2").Value = -100

# Add new row 3
$ws.Range("A3").Value = "27/04/2021"
$ws.Range("B3").Value = "LB0003"
$ws.Range("C3").Value = "LADRILLOS"
$ws.Range("D3").Value = 150
$ws.Range("E3").Value = "BOLSAS"
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = "No"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "Activo"

# Add new row 4
$ws.Range("A4").Value = "18/04/2021"
$ws.Range("B4").Value = "LB0002"
$ws.Range("C4").Value = "ZAPATOS"
$ws.Range("D4").Value = 1500.67
$ws.Range("E4").Value = "UNIDAD"
$ws.Range("F4").Value = 180
$ws.Range("G4").Value = "No"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = "Inactivo"
